$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 299.06489935018368
$ws.Range("B2").Value = 255.58105286261491
$ws.Range("C2").Value = 16.291181932318644
$ws.Range("D2").Value = 32.801547696621327
$ws.Range("E2").Value = 14.630677554772928
$ws.Range("F2").Value = 29.034089744183696
